# Insert 3 new price rows at row 421 on the active sheet, shifting the
# existing rows 421:518 down to 424:521 (dimension grows from A1:T518 to
# A1:T521), then populate the 3 newly inserted rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 421 (each Insert() pushes
# everything at/below the target index down by one row).
$ws.Rows.Item(421).EntireRow.Insert()
$ws.Rows.Item(421).EntireRow.Insert()
$ws.Rows.Item(421).EntireRow.Insert()

# New row 421: Maduro
$ws.Cells.Item(421, 1).Value = 3
$ws.Cells.Item(421, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(421, 3).Value = "Coquimbo"
$ws.Cells.Item(421, 4).Value = 44508
$ws.Cells.Item(421, 5).Value = 5
$ws.Cells.Item(421, 6).Value = "Fruta"
$ws.Cells.Item(421, 7).Value = 100108
$ws.Cells.Item(421, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(421, 9).Value = 100108006
$ws.Cells.Item(421, 10).Value = "Plátano"
$ws.Cells.Item(421, 11).Value = "Sin especificar"
$ws.Cells.Item(421, 12).Value = "Maduro"
$ws.Cells.Item(421, 13).Value = 200
$ws.Cells.Item(421, 14).Value = 17000
$ws.Cells.Item(421, 15).Value = 17000
$ws.Cells.Item(421, 16).Value = 17000
$ws.Cells.Item(421, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(421, 18).Value = "Ecuador"
$ws.Cells.Item(421, 19).Value = 850
$ws.Cells.Item(421, 20).Value = 20

# New row 422: Pintón
$ws.Cells.Item(422, 1).Value = 3
$ws.Cells.Item(422, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(422, 3).Value = "Coquimbo"
$ws.Cells.Item(422, 4).Value = 44508
$ws.Cells.Item(422, 5).Value = 5
$ws.Cells.Item(422, 6).Value = "Fruta"
$ws.Cells.Item(422, 7).Value = 100108
$ws.Cells.Item(422, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(422, 9).Value = 100108006
$ws.Cells.Item(422, 10).Value = "Plátano"
$ws.Cells.Item(422, 11).Value = "Sin especificar"
$ws.Cells.Item(422, 12).Value = "Pintón"
$ws.Cells.Item(422, 13).Value = 320
$ws.Cells.Item(422, 14).Value = 18000
$ws.Cells.Item(422, 15).Value = 18000
$ws.Cells.Item(422, 16).Value = 18000
$ws.Cells.Item(422, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(422, 18).Value = "Ecuador"
$ws.Cells.Item(422, 19).Value = 900
$ws.Cells.Item(422, 20).Value = 20

# New row 423: Primera Pintón
$ws.Cells.Item(423, 1).Value = 3
$ws.Cells.Item(423, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(423, 3).Value = "Coquimbo"
$ws.Cells.Item(423, 4).Value = 44508
$ws.Cells.Item(423, 5).Value = 5
$ws.Cells.Item(423, 6).Value = "Fruta"
$ws.Cells.Item(423, 7).Value = 100108
$ws.Cells.Item(423, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(423, 9).Value = 100108006
$ws.Cells.Item(423, 10).Value = "Plátano"
$ws.Cells.Item(423, 11).Value = "Sin especificar"
$ws.Cells.Item(423, 12).Value = "Primera Pintón"
$ws.Cells.Item(423, 13).Value = 280
$ws.Cells.Item(423, 14).Value = 20000
$ws.Cells.Item(423, 15).Value = 20000
$ws.Cells.Item(423, 16).Value = 20000
$ws.Cells.Item(423, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(423, 18).Value = "Ecuador"
$ws.Cells.Item(423, 19).Value = 1000
$ws.Cells.Item(423, 20).Value = 20
